$wb = $excel.ActiveWorkbook

# --- Sheet1: library_content ---
$ws1 = $wb.Worksheets.Item("library_content")
$ws1.Range("B2").Value = 3

# --- Sheet2: spec ---
$ws2 = $wb.Worksheets.Item("spec")
# Fix the mistranslated French probability labels ("unlikely" must read as
# "improbable", not the opposite meaning "probable").
$ws2.Range("L5").Value = "2 - plutôt improbable"
$ws2.Range("L6").Value = "1 - improbable"

# Widen column L so the longer French label fits without truncation.
# (COM ColumnWidth uses "characters" units with a fixed +5/6 offset vs the
# stored OOXML width, so back that offset out to land exactly on 18.)
$ws2.Columns("L").ColumnWidth = 17.166666666666668

# Activate "spec" first, set its scroll/selection state, then activate
# "library_content" last so it ends up as the selected tab when the file
# is saved (matches the reviewer re-checking the fixed labels then
# returning focus to the first sheet).
$ws2.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 6
$ws2.Range("L6").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("B3").Select() | Out-Null
